# Update the scheduled post timestamps in column F for rows 16-38.
# These rows previously all shared the same datetime
# (2024-04-09 19:22:00 -> Excel serial 45391.806944444441); the refreshed
# PostExecutor / DataHolder run re-scheduled them to
# 2024-04-11 21:05:00 -> Excel serial 45393.878472222219.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = 45393.878472222219

for ($row = 16; $row -le 38; $row++) {
    $ws.Cells.Item($row, 6).Value2 = $newTimestamp
}
